$d = $word.ActiveDocument

# 1. "mix of centers operating" -> "mix of centres operating"
$d.Content.Find.Execute("mix of centers operating", $true, $false, $false, $false, $false,
                         $true, 1, $false, "mix of centres operating", 2) | Out-Null

# 2. "not all centers have" -> "not all centres have"
$d.Content.Find.Execute("not all centers have", $true, $false, $false, $false, $false,
                         $true, 1, $false, "not all centres have", 2) | Out-Null

# 3. "from Centers not having" -> "from Centres not having"
$d.Content.Find.Execute("from Centers not having", $true, $false, $false, $false, $false,
                         $true, 1, $false, "from Centres not having", 2) | Out-Null

# 4. "The Center will decide" -> "The Centre will decide"
$d.Content.Find.Execute("The Center will decide", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The Centre will decide", 2) | Out-Null

# 5. "A center not having made" -> "A centre not having made"
$d.Content.Find.Execute("A center not having made", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A centre not having made", 2) | Out-Null
